$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated Coin / Link / Price / Volume(1h) values row by row
$ws.Range("D2").Value = "61.581.89"
$ws.Range("E2").Value = "  +1.00%  "

$ws.Range("D3").Value = "3.448.39"
$ws.Range("E3").Value = "  +2.14%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.83"
$ws.Range("E5").Value = "  +1.36%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.10"
$ws.Range("E6").Value = "  +6.40%  "

$ws.Range("D7").Value = "3.449.79"
$ws.Range("E7").Value = "  +2.24%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.477"
$ws.Range("E9").Value = "  +1.78%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.63"
$ws.Range("E10").Value = "  +0.07%  "

$ws.Range("E11").Value = "  +3.26%  "

$ws.Range("E12").Value = "  +2.51%  "

$ws.Range("D13").Value = "4.039.89"
$ws.Range("E13").Value = "  +2.23%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.99"
$ws.Range("E14").Value = "  +9.00%  "

$ws.Range("E15").Value = "  -0.88%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.466.59"
$ws.Range("E16").Value = "  +2.72%  "

$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000174"
$ws.Range("E17").Value = "  +1.81%  "

$ws.Range("D18").Value = "61.728.85"
$ws.Range("E18").Value = "  +0.97%  "

$ws.Range("E19").Value = "  +8.96%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.30"
$ws.Range("E20").Value = "  +3.78%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.56"
$ws.Range("E21").Value = "  +2.55%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "389.25"
$ws.Range("E22").Value = "  +3.73%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.565"
$ws.Range("E23").Value = "  +3.39%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.49"
$ws.Range("E24").Value = "  +3.47%  "

$ws.Range("B25").Value = "LEO"
$ws.Range("C25").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.77"
$ws.Range("E25").Value = "  +0.20%  "

$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.995"
$ws.Range("E26").Value = "  -0.65%  "

$ws.Range("E27").Value = "  +0.45%  "

$ws.Range("D28").Value = "3.594.16"
$ws.Range("E28").Value = "  +2.14%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.181"
$ws.Range("E29").Value = "  +2.67%  "

$ws.Range("E30").Value = "  +3.22%  "

$ws.Range("E31").Value = "  +0.13%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.17"
$ws.Range("E32").Value = "  +1.58%  "

$ws.Range("E33").Value = "  +2.39%  "

$ws.Range("E34").Value = "  -11.04%  "

$ws.Range("E35").Value = "  -0.01%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "24.06"
$ws.Range("E36").Value = "  +2.98%  "

$ws.Range("D37").Value = "3.477.78"
$ws.Range("E37").Value = "  +2.39%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.00"
$ws.Range("E38").Value = "  +3.20%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.13"
$ws.Range("E39").Value = "  +0.29%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.55"
$ws.Range("E40").Value = "  +0.56%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "166.67"
$ws.Range("E41").Value = "  +1.26%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0782"
$ws.Range("E42").Value = "  +3.11%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "27.91"
$ws.Range("E43").Value = "  +11.46%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.805"
$ws.Range("E44").Value = "  +3.63%  "

$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.04%  "

$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "42.43"
$ws.Range("E46").Value = "  +1.82%  "

$ws.Range("B47").Value = "Filecoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.49"
$ws.Range("E47").Value = "  +4.14%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.72"
$ws.Range("E48").Value = "  +2.66%  "

$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "2.572.59"
$ws.Range("E49").Value = "  +1.34%  "

$ws.Range("B50").Value = "ONDO"
$ws.Range("C50").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.16"
$ws.Range("E50").Value = "  -1.68%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.93"
$ws.Range("E51").Value = "  +2.47%  "

